# Apply the ETL-refactor data update to the absenteeism worksheet.
# Only the values in rows 2-11 (columns A-G) change; headers and
# formatting/styles stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (row number => A..G values)
$data = @(
    @{ Row = 2;  A = 81670; B = "Lucas Gabriel Cardoso";    C = "Marketing";              D = "Viagem de negócios"; E = 7; F = 45088; G = 4994.29 }
    @{ Row = 3;  A = 56802; B = "Gabriela Araújo";           C = "Vendas";                 D = "Doença";             E = 4; F = 45080; G = 4539.36 }
    @{ Row = 4;  A = 91245; B = "João Gabriel Cavalcanti";   C = "Engenharia";             D = "Doença";             E = 2; F = 45094; G = 7091.92 }
    @{ Row = 5;  A = 40045; B = "Gabrielly Moraes";          C = "Atendimento ao Cliente"; D = "Viagem de negócios"; E = 8; F = 45106; G = 4797.08 }
    @{ Row = 6;  A = 4099;  B = "Sr. Caio Carvalho";         C = "TI";                     D = "Problemas pessoais"; E = 6; F = 45098; G = 4250.25 }
    @{ Row = 7;  A = 72394; B = "Evelyn Souza";              C = "P&D";                    D = "Consulta médica";    E = 8; F = 45089; G = 7302.44 }
    @{ Row = 8;  A = 23992; B = "Brenda Aragão";             C = "Engenharia";             D = "Problemas pessoais"; E = 5; F = 45100; G = 3502.64 }
    @{ Row = 9;  A = 2696;  B = "Isaac da Rosa";             C = "Operações";              D = "Viagem de negócios"; E = 5; F = 45099; G = 12339.88 }
    @{ Row = 10; A = 81177; B = "Maria Vitória Lima";        C = "Atendimento ao Cliente"; D = "Problemas pessoais"; E = 1; F = 45086; G = 4049.69 }
    @{ Row = 11; A = 52802; B = "Julia Silveira";            C = "Marketing";              D = "Viagem de negócios"; E = 4; F = 45085; G = 6459.05 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
}
